# Update Global_M2 "New Zealand_FX" sheet with the latest monthly FX bars.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the (previously provisional) high/close for the existing last row ---
$ws.Range("D210").Value = 1.6356
$ws.Range("F210").Value = 1.6172

# --- Append three new monthly rows (211-213) ---

# Row 211
$ws.Range("A210").Copy()
$ws.Range("A211").PasteSpecial(-4122)
$ws.Range("A211").Value = 45047.33333333334
$ws.Range("B211").Value = "FX_IDC:USDNZD"
$ws.Range("C211").Value = 1.6163
$ws.Range("D211").Value = 1.6702
$ws.Range("E211").Value = 1.5669
$ws.Range("F211").Value = 1.6614
$ws.Range("G211").Value = 0

# Row 212
$ws.Range("A210").Copy()
$ws.Range("A212").PasteSpecial(-4122)
$ws.Range("A212").Value = 45078.33333333334
$ws.Range("B212").Value = "FX_IDC:USDNZD"
$ws.Range("C212").Value = 1.6614
$ws.Range("D212").Value = 1.6685
$ws.Range("E212").Value = 1.6013
$ws.Range("F212").Value = 1.6288
$ws.Range("G212").Value = 0

# Row 213
$ws.Range("A210").Copy()
$ws.Range("A213").PasteSpecial(-4122)
$ws.Range("A213").Value = 45110.33333333334
$ws.Range("B213").Value = "FX_IDC:USDNZD"
$ws.Range("C213").Value = 1.6288
$ws.Range("D213").Value = 1.634
$ws.Range("E213").Value = 1.6082
$ws.Range("F213").Value = 1.6101
$ws.Range("G213").Value = 0
